$d = $word.ActiveDocument

# 1) The non-gaussian ranovas were re-run (1000 iterations instead of 10),
#    which changed the chi-square symbol's glyph encoding in the two
#    "header1" table header rows ("Χ2" column) from χ (U+03C7) to the
#    mis-encoded Ï‡ (U+00CF U+2021) sequence found in the regenerated
#    table export.
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute([string][char]0x03C7, $false, $false, $false, $false, $false, `
                         $true, 1, $false, ([string][char]0x00CF + [string][char]0x2021), 2)

# 2) The regenerated tables also reflowed: the header row height in each
#    of the two affected tables grew from 571 to 637 twips (28.55pt -> 31.85pt).
for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $tbl = $d.Tables.Item($t)
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $row = $tbl.Rows.Item($r)
        if ([Math]::Round($row.Height * 100) -eq 2855) {
            $row.Height = 31.85
        }
    }
}
